# Auto-generated edit script: rewrites the factor-analysis example/
# documentation tables to match the refreshed analysis run referenced
# in the commit "Update examples and documentation".

$wb = $excel.ActiveWorkbook

# --- "Significant Components": refresh per-factor variable-name lists ---
$ws1 = $wb.Worksheets.Item("Significant Components")
$ws1.Range("C2").Value = '[''QEDLESHI'' ''PPUNIT'' ''QHISPC'' ''QEXTRCT'' ''QESL'' ''QNOHLTH'' ''QPOVTY'' ''PERCAP'']'
$ws1.Range("C3").Value = '[''QPOVTY'' ''QSERV'' ''QFHH'' ''QBLACK'' ''QFAM'' ''PERCAP'' ''QRICH'']'
$ws1.Range("C4").Value = '[''QAGEDEP'' ''QFEMALE'' ''QFEMLBR'']'
$ws1.Range("C5").Value = '[''QRENTER'' ''QAGEDEP'' ''MEDAGE'' ''QSSBEN'']'

# --- "Loading Factors": refreshed loadings, rows re-sorted by factor ---
$ws2 = $wb.Worksheets.Item("Loading Factors")
$lf = New-Object 'object[,]' 20,6
$lf[0,0] = 'QEDLESHI'
$lf[0,1] = 0.8692554782357783
$lf[0,2] = 0.2820098258644139
$lf[0,3] = -0.09611805589828241
$lf[0,4] = -0.04435359344350712
$lf[0,5] = 0.01766608403865631
$lf[1,0] = 'PPUNIT'
$lf[1,1] = 0.5203431740155542
$lf[1,2] = 0.1723556441346736
$lf[1,3] = -0.008531185261721667
$lf[1,4] = -0.04118569706715453
$lf[1,5] = -0.6351977716203236
$lf[2,0] = 'QHISPC'
$lf[2,1] = 0.8803370957544391
$lf[2,2] = 0.1202941239550452
$lf[2,3] = -0.1034740118703461
$lf[2,4] = -0.2055870767799969
$lf[2,5] = -0.1433395973189483
$lf[3,0] = 'QEXTRCT'
$lf[3,1] = 0.7518642473093413
$lf[3,2] = 0.1080633735595425
$lf[3,3] = -0.1841536326290315
$lf[3,4] = -0.04880562998113309
$lf[3,5] = 0.02811185768577718
$lf[4,0] = 'QESL'
$lf[4,1] = 0.866987049019839
$lf[4,2] = 0.1055887591660848
$lf[4,3] = -0.1218582046621672
$lf[4,4] = -0.1645406170571051
$lf[4,5] = 0.124620247327293
$lf[5,0] = 'QNOHLTH'
$lf[5,1] = 0.7427426207000754
$lf[5,2] = 0.3825258297491037
$lf[5,3] = -0.09377154715997181
$lf[5,4] = -0.1242320393971812
$lf[5,5] = 0.1047037140321911
$lf[6,0] = 'QPOVTY'
$lf[6,1] = 0.4926703471130989
$lf[6,2] = 0.4936971007074465
$lf[6,3] = 0.01006171172543471
$lf[6,4] = -0.1399809107396388
$lf[6,5] = 0.3593031970426002
$lf[7,0] = 'QSERV'
$lf[7,1] = 0.3709861026223996
$lf[7,2] = 0.5379072882887543
$lf[7,3] = 0.00575919881979431
$lf[7,4] = -0.0875776964683645
$lf[7,5] = 0.1559839714259979
$lf[8,0] = 'QFHH'
$lf[8,1] = 0.2285252628627771
$lf[8,2] = 0.7131346598683198
$lf[8,3] = 0.2214646500776639
$lf[8,4] = -0.07783278026777046
$lf[8,5] = -0.04601151877673739
$lf[9,0] = 'QBLACK'
$lf[9,1] = -0.2749361524699589
$lf[9,2] = 0.704648556134802
$lf[9,3] = 0.05019711367311201
$lf[9,4] = 0.1307361797650962
$lf[9,5] = 0.1670202060833451
$lf[10,0] = 'QFAM'
$lf[10,1] = 0.1243576373666156
$lf[10,2] = 0.6629400806657814
$lf[10,3] = 0.07601969973533645
$lf[10,4] = -0.118412562716169
$lf[10,5] = 0.2379118982834663
$lf[11,0] = 'PERCAP'
$lf[11,1] = 0.5037823274519664
$lf[11,2] = 0.706143335809101
$lf[11,3] = -0.04688443578342617
$lf[11,4] = -0.1072480917650937
$lf[11,5] = -0.1020625160227645
$lf[12,0] = 'QRICH'
$lf[12,1] = 0.403877884430713
$lf[12,2] = 0.6463639841834479
$lf[12,3] = -0.06158260947039726
$lf[12,4] = -0.1194534714731934
$lf[12,5] = 0.002273120718433929
$lf[13,0] = 'QNOAUTO'
$lf[13,1] = 0.189130410502966
$lf[13,2] = 0.3968119687604864
$lf[13,3] = -0.001312813302421819
$lf[13,4] = 0.08291571109167226
$lf[13,5] = 0.5521493329099174
$lf[14,0] = 'QRENTER'
$lf[14,1] = 0.17272274185355
$lf[14,2] = 0.3584593782385577
$lf[14,3] = -0.05642809085152523
$lf[14,4] = -0.4605499124871737
$lf[14,5] = 0.657969995084253
$lf[15,0] = 'QAGEDEP'
$lf[15,1] = -0.1133192111059629
$lf[15,2] = -0.0673948075316152
$lf[15,3] = 0.7290403985390425
$lf[15,4] = 0.4777003030203259
$lf[15,5] = 0.06612983556966748
$lf[16,0] = 'QFEMALE'
$lf[16,1] = -0.1201784064898982
$lf[16,2] = 0.08442001436012979
$lf[16,3] = 0.9476421078584786
$lf[16,4] = 0.01306748466847658
$lf[16,5] = -0.03709586634662605
$lf[17,0] = 'QFEMLBR'
$lf[17,1] = -0.4233895396708394
$lf[17,2] = 0.2565876918738777
$lf[17,3] = 0.5549136198598398
$lf[17,4] = -0.009510213121744099
$lf[17,5] = -0.03382138775670606
$lf[18,0] = 'MEDAGE'
$lf[18,1] = -0.3065470996556154
$lf[18,2] = -0.3702019272504437
$lf[18,3] = 0.04864970456986826
$lf[18,4] = 0.6451312876761621
$lf[18,5] = 0.09075703990682774
$lf[19,0] = 'QSSBEN'
$lf[19,1] = -0.08069646329989939
$lf[19,2] = 0.07089796186361795
$lf[19,3] = 0.1345297509011206
$lf[19,4] = 0.8476251258068529
$lf[19,5] = -0.1136918281150361
$ws2.Range("A2:F21").Value = $lf

# --- "All Refactor Variances": refreshed SS loadings / variance stats ---
$ws3 = $wb.Worksheets.Item("All Refactor Variances")
$rv = New-Object 'object[,]' 4,12
$rv[0,0] = 4.794886974007128
$rv[0,1] = 2.924553900504404
$rv[0,2] = 2.22402164593982
$rv[0,3] = 1.964387356238504
$rv[0,4] = 1.939969486264947
$rv[0,5] = 1.919869526961092
$rv[0,6] = 0.6695949908256502
$rv[0,7] = 4.983165037946606
$rv[0,8] = 3.69248034038259
$rv[0,9] = 1.902149222502244
$rv[0,10] = 1.763128723633083
$rv[0,11] = 1.4679336985806
$rv[1,0] = 0.1775884064447084
$rv[1,1] = 0.1083168111297927
$rv[1,2] = 0.0823711720718452
$rv[1,3] = 0.07275508726809274
$rv[1,4] = 0.07185072171351656
$rv[1,5] = 0.07110627877633673
$rv[1,6] = 0.02479981447502408
$rv[1,7] = 0.2491582518973303
$rv[1,8] = 0.1846240170191295
$rv[1,9] = 0.09510746112511219
$rv[1,10] = 0.08815643618165415
$rv[1,11] = 0.07339668492903
$rv[2,0] = 0.1775884064447084
$rv[2,1] = 0.2859052175745012
$rv[2,2] = 0.3682763896463463
$rv[2,3] = 0.4410314769144391
$rv[2,4] = 0.5128821986279557
$rv[2,5] = 0.5839884774042924
$rv[2,6] = 0.6087882918793165
$rv[2,7] = 0.2491582518973303
$rv[2,8] = 0.4337822689164598
$rv[2,9] = 0.528889730041572
$rv[2,10] = 0.6170461662232262
$rv[2,11] = 0.6904428511522561
$rv[3,0] = 0.2917079858689411
$rv[3,1] = 0.1779219682353303
$rv[3,2] = 0.1353034760533373
$rv[3,3] = 0.1195080264167028
$rv[3,4] = 0.1180225090921425
$rv[3,5] = 0.1167996818020812
$rv[3,6] = 0.04073635253146473
$rv[3,7] = 0.3608673063694102
$rv[3,8] = 0.2673994186644367
$rv[3,9] = 0.137748491372444
$rv[3,10] = 0.1276810036261987
$rv[3,11] = 0.1063037799675104
$ws3.Range("B2:M5").Value = $rv

# --- "Final Variances": refreshed SS loadings / variance stats ---
$ws4 = $wb.Worksheets.Item("Final Variances")
$fv = New-Object 'object[,]' 4,5
$fv[0,0] = 4.983165037946606
$fv[0,1] = 3.69248034038259
$fv[0,2] = 1.902149222502244
$fv[0,3] = 1.763128723633083
$fv[0,4] = 1.4679336985806
$fv[1,0] = 0.2491582518973303
$fv[1,1] = 0.1846240170191295
$fv[1,2] = 0.09510746112511219
$fv[1,3] = 0.08815643618165415
$fv[1,4] = 0.07339668492903
$fv[2,0] = 0.2491582518973303
$fv[2,1] = 0.4337822689164598
$fv[2,2] = 0.528889730041572
$fv[2,3] = 0.6170461662232262
$fv[2,4] = 0.6904428511522561
$fv[3,0] = 0.3608673063694102
$fv[3,1] = 0.2673994186644367
$fv[3,2] = 0.137748491372444
$fv[3,3] = 0.1276810036261987
$fv[3,4] = 0.1063037799675104
$ws4.Range("B2:F5").Value = $fv

# --- "Included and Excluded": refresh the include-list ordering ---
$ws5 = $wb.Worksheets.Item("Included and Excluded")
$ws5.Range("B2").Value = '[[''QEDLESHI'', ''PPUNIT'', ''QHISPC'', ''QEXTRCT'', ''QESL'', ''QNOHLTH'', ''QPOVTY'', ''PERCAP'', ''QSERV'', ''QFHH'', ''QBLACK'', ''QFAM'', ''QRICH'', ''QAGEDEP'', ''QFEMALE'', ''QFEMLBR'', ''QRENTER'', ''MEDAGE'', ''QSSBEN'', ''QNOAUTO'']]'

